# Added New Mac-Address and Document Types
# Appends 5 new device_master_h rows (ids 3000176-3000180) for "32"-suffixed
# devices, mirroring the existing row pattern, and reserves a handful of
# additional formatted-but-empty rows below them (as left behind by the
# original author's worksheet navigation / fill pattern).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlLeft alignment constant used for column H formatting (matches existing rows)
$xlLeft = -4131

$newRows = @(
    @{ Row = 157; Id = 3000176; Name = "Finger Print Scanner 32"; Mac = "80-75-40-E8-CA-24"; Serial = "BS563Q2230824"; DSpec = 165 },
    @{ Row = 158; Id = 3000177; Name = "IRIS Scanner 32";          Mac = "0E-1A-14-4A-6D-3A"; Serial = "BS563Q2230825"; DSpec = 327 },
    @{ Row = 159; Id = 3000178; Name = "Web Camera 32";            Mac = "65-13-7F-0F-F7-53"; Serial = "BS563Q2230826"; DSpec = 736 },
    @{ Row = 160; Id = 3000179; Name = "Document Scanner 32";      Mac = "73-C4-DE-8E-C9-8D"; Serial = "BS563Q2230827"; DSpec = 801 },
    @{ Row = 161; Id = 3000180; Name = "Printer 32";                Mac = "EC-74-AB-E0-0F-38"; Serial = "BS563Q2230828"; DSpec = 920 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Id        # A: id
    $ws.Cells.Item($row, 2).Value = $r.Name       # B: name
    $ws.Cells.Item($row, 3).Value = $r.Mac        # C: mac_address
    $ws.Cells.Item($row, 4).Value = $r.Serial     # D: serial_num
    $ws.Cells.Item($row, 6).Value = $r.DSpec      # F: dspec_id
    $ws.Cells.Item($row, 7).Value = "eng"         # G: lang_code
    $ws.Cells.Item($row, 8).Value = $true         # H: is_active
    $ws.Cells.Item($row, 8).HorizontalAlignment = $xlLeft
    $ws.Cells.Item($row, 9).Value = "superadmin"  # I: cr_by
    $ws.Cells.Item($row, 10).Value = "now()"      # J: cr_dtimes
    $ws.Cells.Item($row, 11).Value = "now()"      # K: eff_dtimes
}

# Trailing rows left with only the column-H left-aligned style applied
for ($row = 162; $row -le 166; $row++) {
    $ws.Cells.Item($row, 8).HorizontalAlignment = $xlLeft
}

# Scroll/select to mirror the author's final on-screen position
$excel.Goto($ws.Range("A154"), $true) | Out-Null
$ws.Range("E159").Select() | Out-Null
